$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "217.43",
# "4.510", "1.000") are preserved exactly as text and are not auto-converted to
# numbers (which would drop significant trailing zeros / change the stored type).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.062.50"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.651.27"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "217.43"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "0.5267"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.2595"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "0.06308"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "20.31"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "4.510"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "1.668.55"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "1.878.27"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "0.5482"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "26.072.72"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "4.565"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "190.42"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "6.019"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "142.59"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "0.1237"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "7.227"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "16.21"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "1.431"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "0.05816"
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("D31").Value = "1.271"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "3.546"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "3.258"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "1.590"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "2.793"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "0.9424"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "0.5762"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.8501"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "104.85"
$ws.Range("E41").Value = "  +4.45%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "5.723"
$ws.Range("E43").Value = "  -4.93%  "
$ws.Range("D44").Value = "1.033.05"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "1.792.92"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "57.12"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "0.4325"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "0.05140"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "7.796"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").Value = "1.452"
$ws.Range("E51").Value = "  -0.50%  "
